$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.917.86"
$ws.Range("E2").Value = "  +2.20%  "
$ws.Range("D3").Value = "1.706.82"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'312.78"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.27%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D8").Value = "'49.40"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.95%  "
$ws.Range("D9").Value = "'0.3441"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").Value = "'1.224"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.53%  "
$ws.Range("D11").Value = "'0.07533"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'1.000"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "'21.19"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.49%  "
$ws.Range("D14").Value = "'6.345"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.29%  "
$ws.Range("D15").Value = "'7.052"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +5.03%  "
$ws.Range("D16").Value = "1.709.57"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").Value = "'0.00001133"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").Value = "'0.06720"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").Value = "'0.9990"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "'84.08"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.38%  "
$ws.Range("D21").Value = "'17.37"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.79%  "
$ws.Range("D22").Value = "'6.377"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.62%  "
$ws.Range("D23").Value = "'13.27"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +8.82%  "
$ws.Range("D24").Value = "24.915.89"
$ws.Range("E24").Value = "  +2.32%  "
$ws.Range("D25").Value = "'2.450"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  +5.75%  "
$ws.Range("D27").Value = "'20.38"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.78%  "
$ws.Range("D28").Value = "'149.84"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("D29").Value = "'132.86"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.95%  "
$ws.Range("B30").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C30").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D30").Value = "1.896.45"
$ws.Range("E30").Value = "  +1.90%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.251"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +29.23%  "
$ws.Range("D32").Value = "'6.819"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +8.68%  "
$ws.Range("D33").Value = "'4.230"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.46%  "
$ws.Range("D34").Value = "'13.79"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +13.00%  "
$ws.Range("D35").Value = "'1.779"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.71%  "
$ws.Range("D36").Value = "'0.08797"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.33%  "
$ws.Range("D37").Value = "'5.616"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.78%  "
$ws.Range("D38").Value = "'0.06663"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.46%  "
$ws.Range("D39").Value = "'9.184"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.33%  "
$ws.Range("D40").Value = "'0.02420"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.51%  "
$ws.Range("D41").Value = "'0.2251"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +7.99%  "
$ws.Range("D42").Value = "'1.273"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.87%  "
$ws.Range("D43").Value = "'0.6478"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +6.47%  "
$ws.Range("D44").Value = "'0.9992"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").Value = "'13.84"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +6.88%  "
$ws.Range("D46").Value = "'0.6166"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.94%  "
$ws.Range("D47").Value = "'3.837"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.44%  "
$ws.Range("D48").Value = "'2.120"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.24%  "
$ws.Range("D49").Value = "'129.31"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("D50").Value = "'0.07325"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.59%  "
$ws.Range("D51").Value = "'80.30"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.21%  "
